$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected the hours logged on 01/02/2018 for Giovanni (row 39): 1h -> 8h
$ws.Range("D39").Value = 0.33333333333333331

# Added the missing entry for Mirko on the same day (hierarchy management)
$ws.Range("A40").Value = 43132
$ws.Range("B40").Value = "Mirko"
$ws.Range("C40").Value = "gestione gerarchia"
$ws.Range("D40").Value = 0.4375
